# Update gh-pages to output generated at 456a3b4
# Applies numeric "F" (and a couple of "G") column refreshes across the
# four sheets of 北京-漫展信息.xlsx, updates the 本地生活 sold-out status
# text, and removes the duplicate last row on 全部类型.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 210
$ws.Cells.Item(6, 6).Value = 7
$ws.Cells.Item(7, 6).Value = 757
$ws.Cells.Item(9, 6).Value = 9810
$ws.Cells.Item(11, 6).Value = 2656
$ws.Cells.Item(11, 7).Value = 49
$ws.Cells.Item(13, 6).Value = 2392
$ws.Cells.Item(14, 6).Value = 2661
$ws.Cells.Item(16, 6).Value = 278
$ws.Cells.Item(17, 6).Value = 2092
$ws.Cells.Item(19, 6).Value = 82
$ws.Cells.Item(20, 6).Value = 368
$ws.Cells.Item(23, 6).Value = 301
$ws.Cells.Item(25, 6).Value = 154
$ws.Cells.Item(26, 6).Value = 599
$ws.Cells.Item(27, 6).Value = 1291
$ws.Cells.Item(32, 6).Value = 1684
$ws.Cells.Item(33, 6).Value = 2824
$ws.Cells.Item(35, 6).Value = 997
$ws.Cells.Item(36, 6).Value = 360
$ws.Cells.Item(39, 6).Value = 55
$ws.Cells.Item(43, 6).Value = 27

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(14, 6).Value = 157

# ---------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 725
$ws.Cells.Item(3, 6).Value = 954
$ws.Cells.Item(5, 6).Value = 1762
$ws.Cells.Item(5, 7).Value = "已售罄"

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 210
$ws.Cells.Item(3, 6).Value = 725
$ws.Cells.Item(4, 6).Value = 954
$ws.Cells.Item(10, 6).Value = 757
$ws.Cells.Item(12, 6).Value = 9810
$ws.Cells.Item(16, 6).Value = 2656
$ws.Cells.Item(16, 7).Value = 49
$ws.Cells.Item(18, 6).Value = 2392
$ws.Cells.Item(19, 6).Value = 2661
$ws.Cells.Item(20, 6).Value = 278
$ws.Cells.Item(21, 6).Value = 2092
$ws.Cells.Item(23, 6).Value = 82
$ws.Cells.Item(24, 6).Value = 368
$ws.Cells.Item(26, 6).Value = 301
$ws.Cells.Item(28, 6).Value = 154
$ws.Cells.Item(29, 6).Value = 599
$ws.Cells.Item(30, 6).Value = 1291
$ws.Cells.Item(35, 6).Value = 1684
$ws.Cells.Item(37, 6).Value = 2824
$ws.Cells.Item(38, 6).Value = 997
$ws.Cells.Item(41, 6).Value = 360
$ws.Cells.Item(48, 6).Value = 27
$ws.Cells.Item(49, 6).Value = 157

# Row 50 on "全部类型" duplicated row 49 verbatim and was dropped entirely
# (used range shrinks from A1:I50 to A1:I49).
$ws.Rows.Item(50).Delete()
